$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FA2025")

# Replace the generic "Presentation" placeholder with the specific
# presentation group assignments for the three presentation days.
$ws.Range("C30").Value = "Presentation: Shuhang, Catherine, Shreezal"
$ws.Range("C31").Value = "Presentation: Alisha, Suchil, Matt"
$ws.Range("C32").Value = "Presentation: Jayna, Trent, Mehran"

# Move the active selection to reflect where editing left off.
$ws.Range("C33").Select()
